$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: clone row 239 (style/shape donor) into each new row 240-246
$ws.Range("A239:T239").Copy($ws.Range("A240"))
$ws.Range("A239:T239").Copy($ws.Range("A241"))
$ws.Range("A239:T239").Copy($ws.Range("A242"))
$ws.Range("A239:T239").Copy($ws.Range("A243"))
$ws.Range("A239:T239").Copy($ws.Range("A244"))
$ws.Range("A239:T239").Copy($ws.Range("A245"))
$ws.Range("A239:T239").Copy($ws.Range("A246"))

# Step 2: overwrite the cloned cells with the real data for each new row
# Row 240
$ws.Range("A240").Value = 45820
$ws.Range("B240").Value = "Flowering"
$ws.Range("C240").Value = "Large"
$ws.Range("D240").Value = 68
$ws.Range("E240").Value = 82
$ws.Range("G240").Value = 0.68
$ws.Range("H240").Value = 0.1
$ws.Range("I240").Value = "No"
$ws.Range("J240").Value = 2
$ws.Range("K240").Value = "Neutral"
$ws.Range("L240").Value = 6
$ws.Range("M240").Value = 0.67
$ws.Range("N240").Value = 68
$ws.Range("O240").Value = 30.02
$ws.Range("P240").Value = 9
$ws.Range("Q240").Value = 0.65
$ws.Range("R240").Value = 8.1
$ws.Range("S240").Value = 85
$ws.Range("T240").Value = 32
$ws.Range("F240").Formula = "=ABS(D240-E240)"

# Row 241
$ws.Range("A241").Value = 45820
$ws.Range("B241").Value = "Nonflowering"
$ws.Range("C241").Value = "Medium"
$ws.Range("D241").Value = 68
$ws.Range("E241").Value = 82
$ws.Range("G241").Value = 0.68
$ws.Range("H241").Value = 0.2
$ws.Range("I241").Value = "No"
$ws.Range("J241").Value = 3
$ws.Range("K241").Value = "Neutral"
$ws.Range("L241").Value = 6
$ws.Range("M241").Value = 0.67
$ws.Range("N241").Value = 68
$ws.Range("O241").Value = 30.02
$ws.Range("P241").Value = 9
$ws.Range("Q241").Value = 0.65
$ws.Range("R241").Value = 8.1
$ws.Range("S241").Value = 85
$ws.Range("T241").Value = 32
$ws.Range("F241").Formula = "=ABS(D241-E241)"

# Row 242
$ws.Range("A242").Value = 45820
$ws.Range("B242").Value = "Nonflowering"
$ws.Range("C242").Value = "Small"
$ws.Range("D242").Value = 68
$ws.Range("E242").Value = 82
$ws.Range("G242").Value = 0.68
$ws.Range("H242").Value = 0.2
$ws.Range("I242").Value = "No"
$ws.Range("J242").Value = 3
$ws.Range("K242").Value = "Dark"
$ws.Range("L242").Value = 6
$ws.Range("M242").Value = 0.67
$ws.Range("N242").Value = 68
$ws.Range("O242").Value = 30.02
$ws.Range("P242").Value = 9
$ws.Range("Q242").Value = 0.65
$ws.Range("R242").Value = 8.1
$ws.Range("S242").Value = 85
$ws.Range("T242").Value = 32
$ws.Range("F242").Formula = "=ABS(D242-E242)"

# Row 243
$ws.Range("A243").Value = 45820
$ws.Range("B243").Value = "Nonflowering"
$ws.Range("C243").Value = "Medium"
$ws.Range("D243").Value = 68
$ws.Range("E243").Value = 82
$ws.Range("G243").Value = 0.68
$ws.Range("H243").Value = 0.25
$ws.Range("I243").Value = "No"
$ws.Range("J243").Value = 3
$ws.Range("K243").Value = "Neutral"
$ws.Range("L243").Value = 6
$ws.Range("M243").Value = 0.67
$ws.Range("N243").Value = 68
$ws.Range("O243").Value = 30.02
$ws.Range("P243").Value = 9
$ws.Range("Q243").Value = 0.65
$ws.Range("R243").Value = 8.1
$ws.Range("S243").Value = 85
$ws.Range("T243").Value = 32
$ws.Range("F243").Formula = "=ABS(D243-E243)"

# Row 244
$ws.Range("A244").Value = 45820
$ws.Range("B244").Value = "Nonflowering"
$ws.Range("C244").Value = "Medium"
$ws.Range("D244").Value = 68
$ws.Range("E244").Value = 82
$ws.Range("G244").Value = 0.68
$ws.Range("H244").Value = 0.25
$ws.Range("I244").Value = "No"
$ws.Range("J244").Value = 3
$ws.Range("K244").Value = "Bright"
$ws.Range("L244").Value = 6
$ws.Range("M244").Value = 0.67
$ws.Range("N244").Value = 68
$ws.Range("O244").Value = 30.02
$ws.Range("P244").Value = 9
$ws.Range("Q244").Value = 0.65
$ws.Range("R244").Value = 8.1
$ws.Range("S244").Value = 85
$ws.Range("T244").Value = 32
$ws.Range("F244").Formula = "=ABS(D244-E244)"

# Row 245
$ws.Range("A245").Value = 45820
$ws.Range("B245").Value = "Nonflowering"
$ws.Range("C245").Value = "Large"
$ws.Range("D245").Value = 68
$ws.Range("E245").Value = 82
$ws.Range("G245").Value = 0.68
$ws.Range("H245").Value = 0.3
$ws.Range("I245").Value = "No"
$ws.Range("J245").Value = 4
$ws.Range("K245").Value = "Bright"
$ws.Range("L245").Value = 6
$ws.Range("M245").Value = 0.67
$ws.Range("N245").Value = 68
$ws.Range("O245").Value = 30.02
$ws.Range("P245").Value = 9
$ws.Range("Q245").Value = 0.65
$ws.Range("R245").Value = 8.1
$ws.Range("S245").Value = 85
$ws.Range("T245").Value = 32
$ws.Range("F245").Formula = "=ABS(D245-E245)"

# Row 246
$ws.Range("A246").Value = 45820
$ws.Range("B246").Value = "Tree"
$ws.Range("C246").Value = "Medium"
$ws.Range("D246").Value = 68
$ws.Range("E246").Value = 82
$ws.Range("G246").Value = 0.68
$ws.Range("H246").Value = 1.25
$ws.Range("I246").Value = "No"
$ws.Range("J246").Value = 1
$ws.Range("K246").Value = "Neutral"
$ws.Range("L246").Value = 6
$ws.Range("M246").Value = 0.67
$ws.Range("N246").Value = 68
$ws.Range("O246").Value = 30.02
$ws.Range("P246").Value = 9
$ws.Range("Q246").Value = 0.65
$ws.Range("R246").Value = 8.1
$ws.Range("S246").Value = 85
$ws.Range("T246").Value = 32
$ws.Range("F246").Formula = "=ABS(D246-E246)"

# Step 3: match the saved selection / active cell from the source edit
$ws.Range("O240:O246").Select()
